$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Title text change
$ws.Range("A1").Value = "BitmexSwap"

# 2. Column A width (best-effort: nearest reachable value to the target stored width 12.125)
$ws.Columns.Item(1).ColumnWidth = 11.428571428571429

# 3. Row 6 updated figures -- E6 changes, B6/F6 recalc automatically via existing formulas
$ws.Range("E6").Value = 3.7

# 4. Row 7 new entries
$ws.Range("B7").Formula = "=E7-E6-D7"
$ws.Range("C7").Value = 2.25
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 1.99

# 5. Turn F2:F7 into one shared formula block (matches diff's t="shared" group)
$ws.Range("F2:F7").Formula = "=B2/(E2-B2)"

# 6. Copy cell formatting so new / recalculated cells pick up the same conditional-style look
$ws.Range("B3").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 7. I5 now averages the full column instead of just the populated rows
$ws.Range("I5").Formula = "=AVERAGE(F2:F100)"

# 8. Selection moves on to the next empty row
$ws.Range("E8").Select()
